$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking price strings
# (e.g. "1.009", "217.92") are not auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '26.300.15'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '1.678.82'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").Value = '217.92'
$ws.Range("E5").Value = '  +0.73%  '
$ws.Range("D6").Value = '0.5263'
$ws.Range("E6").Value = '  +3.18%  '
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").Value = '0.2689'
$ws.Range("E8").Value = '  +2.36%  '
$ws.Range("D9").Value = '0.06460'
$ws.Range("E9").Value = '  +0.88%  '
$ws.Range("D10").Value = '21.93'
$ws.Range("E10").Value = '  +1.30%  '
$ws.Range("D11").Value = '0.07504'
$ws.Range("E11").Value = '  +1.12%  '
$ws.Range("D12").Value = '1.707.21'
$ws.Range("E12").Value = '  +2.14%  '
$ws.Range("D13").Value = '4.512'
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("D14").Value = '0.5780'
$ws.Range("E14").Value = '  -0.33%  '
$ws.Range("D15").Value = '0.000008502'
$ws.Range("E15").Value = '  -0.06%  '
$ws.Range("D16").Value = '64.77'
$ws.Range("E16").Value = '  +0.80%  '
$ws.Range("D17").Value = '26.339.30'
$ws.Range("E17").Value = '  +0.86%  '
$ws.Range("D18").Value = '4.920'
$ws.Range("E18").Value = '  +0.57%  '
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("D20").Value = '10.87'
$ws.Range("E20").Value = '  +1.21%  '
$ws.Range("D21").Value = '189.57'
$ws.Range("E21").Value = '  +0.61%  '
$ws.Range("D22").Value = '6.196'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").Value = '144.84'
$ws.Range("E24").Value = '  -0.49%  '
$ws.Range("D25").Value = '7.763'
$ws.Range("E25").Value = '  +2.08%  '
$ws.Range("D26").Value = '0.1255'
$ws.Range("E26").Value = '  +5.63%  '
$ws.Range("D27").Value = '15.79'
$ws.Range("E27").Value = '  +1.28%  '
$ws.Range("D28").Value = '0.06526'
$ws.Range("E28").Value = '  +0.66%  '
$ws.Range("E29").Value = '  +4.21%  '
$ws.Range("D30").Value = '1.325'
$ws.Range("E30").Value = '  +0.79%  '
$ws.Range("D31").Value = '3.591'
$ws.Range("E31").Value = '  +1.92%  '
$ws.Range("D32").Value = '3.585'
$ws.Range("E32").Value = '  +2.31%  '
$ws.Range("D33").Value = '1.659'
$ws.Range("E33").Value = '  +2.00%  '
$ws.Range("D34").Value = '1.028'
$ws.Range("E34").Value = '  +1.13%  '
$ws.Range("D35").Value = '0.6204'
$ws.Range("E35").Value = '  +2.59%  '
$ws.Range("D36").Value = '2.406'
$ws.Range("E36").Value = '  +1.70%  '
$ws.Range("E37").Value = '  +2.03%  '
$ws.Range("D38").Value = '6.284'
$ws.Range("E38").Value = '  +1.43%  '
$ws.Range("D39").Value = '1.117.10'
$ws.Range("E39").Value = '  +3.97%  '
$ws.Range("D40").Value = '0.01621'
$ws.Range("E40").Value = '  +0.83%  '
$ws.Range("D41").Value = '0.8731'
$ws.Range("E41").Value = '  +1.75%  '
$ws.Range("E42").Value = '  +0.69%  '
$ws.Range("D43").Value = '100.52'
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '1.830.27'
$ws.Range("E44").Value = '  +0.95%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '56.93'
$ws.Range("E45").Value = '  +1.51%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.00000000106'
$ws.Range("E46").Value = '  -7.34%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").Value = '1.007'
$ws.Range("E47").Value = '  +0.35%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '8.153'
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("D49").Value = '0.05269'
$ws.Range("E49").Value = '  +1.19%  '
$ws.Range("D50").Value = '0.4294'
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("D51").Value = '6.079'
$ws.Range("E51").Value = '  +2.53%  '

# Restore default cell style on the price column so formatting matches
# the original workbook (text values, no explicit number format).
$priceRange.Style = "Normal"
